$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 51
$newRow = $ws.Range("A" + $row + ":C" + $row)

# "2025-10-01" looks like a date, so a plain .Value assignment would get
# auto-detected and silently converted to a date serial number, same as
# typing it into a General-formatted cell in real Excel. Force Text format
# on the new row first so the literal strings are preserved verbatim, then
# clear the formatting back off so these cells don't end up with a style
# index that the rest of the sheet's cells (style-less) don't have.
$newRow.NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-10-01"
$ws.Cells.Item($row, 2).Value = "15:19:39"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,723.7537"
$newRow.ClearFormats()
